# Update column F (dSF) values for specific rows on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -9
$ws.Range("F4").Value = -13
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = -9
